# Updates loading_percent values for rows 2-25, columns C-N (commit: "case with 380 kV done")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one row of data for columns C..N (12 columns), for sheet rows 2..25 (24 rows)
$data = @(
  @(6.545057605471523, 4.661537461599443, 12.26468098392145, 56.38140612707684, 3.797787779480617, 0, 43.04087415747635, 11.22231888176117, 27.64467548261242, 10.56920454279153, 0, 21.32664787442998),
  @(6.5267545040153, 4.651014653266074, 12.27464296581825, 56.41359703053492, 3.801942545423381, 0, 43.08750460884202, 11.24756180758913, 27.40371910342969, 10.59012472440471, 0, 21.40816327839199),
  @(6.516341489421612, 4.644726800730176, 12.28243321283558, 56.44721778523787, 3.804625043249062, 0, 43.12618634493339, 11.26439841427812, 27.26136178966464, 10.60411964272775, 0, 21.46032328296774),
  @(6.512308328545109, 4.642209016122701, 12.28602904537922, 56.4643918581226, 3.805751367955448, 0, 43.14446828694417, 11.27159599162501, 27.20481044794085, 10.61011217301629, 0, 21.48211147090877),
  @(6.511651399212753, 4.641793660959926, 12.28665158346781, 56.46745303552255, 3.805940401418025, 0, 43.14765585519477, 11.27281147600191, 27.19550988417102, 10.61112472268948, 0, 21.48576161233748),
  @(6.516286242090123, 4.644692663207821, 12.28248000139093, 56.44743535304995, 3.804640098723071, 0, 43.12642271505108, 11.26449412045022, 27.26059313548969, 10.60419928750616, 0, 21.46061496694973),
  @(6.538577117988623, 4.657873862835729, 12.26776879819972, 56.38962389654733, 3.799193137592886, 0, 43.05486154797011, 11.23074521801831, 27.56046690842269, 10.57617938242119, 0, 21.35431793965384),
  @(6.588710760320669, 4.685051976545656, 12.25217924736337, 56.38660427187746, 3.789548780642032, 0, 42.99461835962765, 11.17516472946253, 28.19035608272177, 10.53034057260008, 0, 21.1625040146238),
  @(6.629281381168896, 4.70577122035801, 12.2487798112409, 56.45212135739484, 3.783087043169787, 0, 42.99959664411205, 11.1407788182739, 28.67510892826655, 10.50219480073777, 0, 21.03157426326242),
  @(6.648508899056186, 4.715347262762391, 12.24897457961468, 56.4966968962224, 3.780281149543867, 0, 43.01262125891291, 11.12653332010594, 28.89962129392635, 10.49058742920866, 0, 20.97415027230877),
  @(6.655897163141986, 4.718994002534502, 12.24929789852083, 56.51570226169178, 3.779237703264413, 0, 43.01910423468225, 11.12133957041184, 28.98514542779465, 10.48636369522084, 0, 20.95271018531318),
  @(6.654301260705835, 4.718207719009524, 12.24921718047813, 56.51151457830193, 3.779461581266573, 0, 43.01763899113843, 11.1224492130343, 28.96670485677181, 10.48726571912368, 0, 20.95731415437044),
  @(6.649114609071037, 4.71564687817014, 12.24899618259327, 56.49821787803668, 3.780194922815002, 0, 43.01312351974551, 11.12610200549948, 28.90664766479684, 10.49023650014611, 0, 20.97238028022518),
  @(6.645951484921911, 4.714080923796233, 12.24889328996853, 56.49035008754743, 3.780646597400739, 0, 43.01055971444286, 11.12836558043701, 28.86992478132112, 10.49207854418213, 0, 20.98164839196829),
  @(6.628040028874557, 4.70514833662234, 12.24880205199074, 56.44950565035153, 3.78327309231216, 0, 42.99896230352518, 11.1417378898964, 28.66051108582832, 10.50297741655641, 0, 21.03536986974604),
  @(6.617246959747231, 4.699706261477708, 12.24919148675711, 56.42823390084467, 3.784918487109465, 0, 42.99460644736818, 11.15029902174716, 28.53301748587471, 10.50996970644736, 0, 21.06887199109862),
  @(6.611111944421884, 4.696590384817748, 12.24957942254662, 56.41738935330321, 3.785877455569716, 0, 42.99311371114995, 11.15535464314879, 28.46006678894306, 10.51410410425916, 0, 21.08834275901077),
  @(6.609047361952558, 4.695537887556177, 12.24973895181621, 56.41395630172462, 3.786204310386546, 0, 42.99278209046239, 11.1570889769094, 28.43543432124271, 10.51552329154645, 0, 21.09496985285705),
  @(6.618388385154654, 4.700284112197688, 12.24913306890962, 56.43035439702592, 3.784742030852306, 0, 42.99496530222659, 11.14937406756643, 28.54655053705429, 10.50921371191489, 0, 21.06528482145444),
  @(6.650635174980736, 4.716398512276389, 12.24905432885424, 56.50206575183181, 3.779979005660417, 0, 43.0144077110436, 11.12502364575583, 28.92427472780086, 10.48935925154428, 0, 20.96794672945042),
  @(6.672333453181321, 4.727049141221304, 12.25045708072568, 56.56132203770859, 3.776977276390304, 0, 43.03615476758986, 11.11027906713302, 29.17405775195246, 10.47738402325813, 0, 20.90610822575601),
  @(6.660696928044923, 4.721354209619622, 12.24957563680823, 56.52856229711912, 3.778569223523826, 0, 43.02371987822331, 11.11804153625184, 29.04049950481738, 10.48368395273157, 0, 20.93895064780285),
  @(6.617872128211727, 4.700022825789522, 12.24915896859553, 56.42939140729165, 3.784821766196068, 0, 42.99479991362097, 11.14979182285491, 28.54043115838077, 10.50955514057408, 0, 21.06690592673712),
  @(6.574478533751543, 4.677563256917405, 12.25497982267642, 56.37555243518763, 3.792047662006469, 0, 43.00229485315689, 11.18906722230059, 28.01583273863062, 10.54176821728009, 0, 21.21262893010058)
)

$nrows = $data.Count
$ncols = $data[0].Count
$arr = New-Object 'object[,]' $nrows, $ncols
for ($i = 0; $i -lt $nrows; $i++) {
  for ($j = 0; $j -lt $ncols; $j++) {
    $arr[$i, $j] = $data[$i][$j]
  }
}

$ws.Range("C2:N25").Value = $arr

Write-Host "Applied loading_percent updates for rows 2-25 (columns C-N)"
